# Update the supervisor list: replace placeholder supervisor names with the
# real supervisor names, fix two priority values, and append four new
# supervisor rows (12-15) with their tid / password / priority.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: replace placeholder supervisor names with real names ---
$ws.Range("A2").Value = "Dr LIU Yang"
$ws.Range("A3").Value = "Dr. ZHANG, Eric Lu"
$ws.Range("A4").Value = "Dr. WAN, Renjie"
$ws.Range("A5").Value = "Dr Yu, Wilson Shih Bun"
$ws.Range("A6").Value = "Prof. CHEUNG, Yiu Ming"
$ws.Range("A7").Value = "Dr. ZHOU, Kaiyang"
$ws.Range("A8").Value = "Prof. Xu, Jianliang"
$ws.Range("A9").Value = "Dr. HAN, Bo"
$ws.Range("A10").Value = "Dr. DAI, Henry Hong Ning"
$ws.Range("A11").Value = "Prof. YUEN, Pong Chi"

# --- Column D: priority corrections for two existing rows ---
$ws.Range("D8").Value = 1
$ws.Range("D11").Value = 2

# --- New rows 12-15: additional supervisors (entered column-by-column,
#     matching how the data was pasted in: password column, then tid
#     column, then name column, then priority column) ---
$ws.Range("C12").Value = "tpw00011"
$ws.Range("C13").Value = "tpw00012"
$ws.Range("C14").Value = "tpw00013"
$ws.Range("C15").Value = "tpw00014"

$ws.Range("B12").Value = "tid00011"
$ws.Range("B13").Value = "tid00012"
$ws.Range("B14").Value = "tid00013"
$ws.Range("B15").Value = "tid00014"

$ws.Range("A12").Value = "Prof. LEUNG,Yiu Wing"
$ws.Range("A13").Value = "Dr. FENG, Jian"
$ws.Range("A14").Value = "Dr. WANG, Juncheng"
$ws.Range("A15").Value = "Dr. YANG, Renchi"

$ws.Range("D12").Value = 2
$ws.Range("D13").Value = 3
$ws.Range("D14").Value = 3
$ws.Range("D15").Value = 3

# --- View state: scroll so row 12 is at top and select D16 (next empty row) ---
$ws.Range("D16").Select()
